$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new parts-list row 30: cover cap of Rituals of Hammam.
$ws.Range("A30").Value = "cover cap of Rituals of Hammam"
$ws.Range("B30").Value = 1
$url = "https://www.amazon.de/RITUALS-Duschschaum-Vorteilspaket-Das-Ritual-Hammam/dp/B0CL4SBDPQ/ref=asc_df_B0CL4SBDPQ/?tag=googshopde-21&linkCode=df0&hvadid=696322103951&hvpos=&hvnetw=g&hvrand=14682189934117811704&hvpone=&hvptwo=&hvqmt=&hvdev=c&hvdvcmdl=&hvlocint=&hvlocphy=9042332&hvtargid=pla-2246768438588&psc=1&mcid=87c684332cef3a189cf8231dafe02054&th=1&psc=1&gad_source=1"
$ws.Range("C30").Value = $url
$ws.Hyperlinks.Add($ws.Range("C30"), $url) | Out-Null
$ws.Range("C30").Style = "Hyperlink"

# Header cell C1 keeps the same visible text "Example Order" but is re-entered
# with a trailing space (matches the new shared string with xml:space="preserve").
$ws.Range("C1").Value = "Example Order "

# Restore the active-cell selection to C1, as in the saved workbook.
$ws.Range("C1").Select() | Out-Null
